# Add a new row (row 60) of logged data to each of the four sheets in the
# workbook, mirroring the structure of the existing row 59 on each sheet.

$wb = $excel.ActiveWorkbook

# Row data for each worksheet, in the same order the sheets appear in the
# workbook: ROW35-FE-LIFTER, ROW35-MID-LIFTER, ROW02-FE-LIFTER, ROW02-MID-LIFTER
$rowsData = @(
    @{
        A = [double]"45754.37155278935"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x6e"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 366
        I = 13
    },
    @{
        A = [double]"45754.22255366898"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x6a"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 362
        I = 14
    },
    @{
        A = [double]"45754.35939628472"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x6e"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 366
        I = 3
    },
    @{
        A = [double]"45754.42365554398"
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x6a"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 362
        I = 3
    }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $data = $rowsData[$i]

    $newRow = 60
    $templateRow = 59

    # Write the new row values, column by column (A..I)
    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I

    # Match the date/time number formatting used in column A of the template row
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($templateRow, 1).NumberFormat
}
